# Update "想去人数" (interest count) figures on the "展览" and "全部类型" sheets
# to match the latest scrape output.
#   Row 6  (银泰百货 exhibition): 129 -> 131
#   Row 10 (第十五届次元之门 expo): 5088 -> 5089

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Cells.Item(6, 6).Value = 131
    $ws.Cells.Item(10, 6).Value = 5089
}
